$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("S:S").Insert()
$ws.Range("S1").Value = "Sub brand"
$ws.AutoFilterMode = $false
$ws.Range("A1:AP32").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Cinema!`$A`$1:`$AP`$32"
    }
}
$ws.Range("S2").Select()
